$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("phen_oncox")

# "EFO and DO update": bump the source_version (column E) for the
# Experimental Factor Ontology (row 4) and Disease Ontology (row 3) rows.
# OncoTree's version in row 2 is left untouched.
$ws.Range("E4").Value = "v3.79.0"
$ws.Range("E3").Value = "v2025-06-27"
